$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3441.4285
$ws.Range("J64").Value = 3513.8462
$ws.Range("L64").Value = 3513.8462
$ws.Range("N64").Value = -4009.8462

$ws.Range("H67").Value = 3441.4285
$ws.Range("J67").Value = 3513.8462
$ws.Range("L67").Value = 3513.8462
$ws.Range("N67").Value = -5229.8462

$ws.Range("H88").Value = 15153542
$ws.Range("I88").Value = 2025.5
$ws.Range("J88").Value = 20835360
$ws.Range("K88").Value = 2025.5
$ws.Range("L88").Value = 20835360
$ws.Range("M88").Value = -1619.5
$ws.Range("N88").Value = -20836172

$ws.Range("H91").Value = 15153542
$ws.Range("I91").Value = 2025.5
$ws.Range("J91").Value = 20835360
$ws.Range("K91").Value = 2025.5
$ws.Range("L91").Value = 20835360
$ws.Range("M91").Value = -621.5
$ws.Range("N91").Value = -20838168

$ws.Range("H118").Value = 1074.1818
$ws.Range("I118").Value = 1190
$ws.Range("J118").Value = 1062.6
$ws.Range("K118").Value = 3570
$ws.Range("L118").Value = 3187.8
$ws.Range("M118").Value = -1913
$ws.Range("N118").Value = -6501.799999999999

$ws.Range("H138").Value = 6668031
$ws.Range("I138").Value = 8547735
$ws.Range("J138").Value = 3627.2727
$ws.Range("K138").Value = 25643205
$ws.Range("L138").Value = 10881.8181
$ws.Range("M138").Value = -25638065
$ws.Range("N138").Value = -21161.8181

$ws.Range("H141").Value = 1709.8628
$ws.Range("I141").Value = 1058.8334
$ws.Range("J141").Value = 3272.3333
$ws.Range("K141").Value = 3176.5002
$ws.Range("L141").Value = 9816.999899999999
$ws.Range("M141").Value = 2003.4998
$ws.Range("N141").Value = -20176.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14379.143
$ws.Range("I32").Value = 14434.548
$ws.Range("J32").Value = 13714.286
$ws.Range("K32").Value = 14434.548
$ws.Range("L32").Value = 13714.286
$ws.Range("M32").Value = -14147.548
$ws.Range("N32").Value = -14288.286

$ws.Range("H61").Value = 1119.8823
$ws.Range("I61").Value = 662.7907
$ws.Range("J61").Value = 3576.75
$ws.Range("K61").Value = 662.7907
$ws.Range("L61").Value = 3576.75
$ws.Range("M61").Value = -450.7907
$ws.Range("N61").Value = -4000.75

$ws.Range("H132").Value = 12518.333
$ws.Range("I132").Value = 17946.143
$ws.Range("J132").Value = 4919.4
$ws.Range("K132").Value = 53838.429
$ws.Range("L132").Value = 14758.2
$ws.Range("M132").Value = -51308.429
$ws.Range("N132").Value = -19818.2

$ws.Range("H136").Value = 1119.8823
$ws.Range("I136").Value = 662.7907
$ws.Range("J136").Value = 3576.75
$ws.Range("K136").Value = 1988.3721
$ws.Range("L136").Value = 10730.25
$ws.Range("M136").Value = 561.6279
$ws.Range("N136").Value = -15830.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 50000
$ws.Range("J57").Value = 50000
$ws.Range("L57").Value = 50000
$ws.Range("N57").Value = -51440

$ws.Range("H107").Value = 950
$ws.Range("I107").Value = 950
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 950
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 13757.5
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 15294.286
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 15294.286
$ws.Range("M36").Value = -2612
$ws.Range("N36").Value = -16070.286

$ws.Range("H40").Value = 13757.5
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 15294.286
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 15294.286
$ws.Range("M40").Value = -2840
$ws.Range("N40").Value = -15614.286

$ws.Range("H58").Value = 860.8
$ws.Range("I58").Value = 672.96
$ws.Range("J58").Value = 1800
$ws.Range("K58").Value = 672.96
$ws.Range("L58").Value = 1800
$ws.Range("M58").Value = -469.96
$ws.Range("N58").Value = -2206

$ws.Range("H107").Value = 720.7143
$ws.Range("I107").Value = 754.86664
$ws.Range("J107").Value = 635.3333
$ws.Range("K107").Value = 754.86664
$ws.Range("L107").Value = 635.3333
$ws.Range("M107").Value = 1165.13336
$ws.Range("N107").Value = -4475.3333

$ws.Range("H136").Value = 860.8
$ws.Range("I136").Value = 672.96
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 2018.88
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = 531.1199999999999
$ws.Range("N136").Value = -10500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 17200
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 17200
$ws.Range("K59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").Value = -52680

$ws.Range("H63").Value = 6840
$ws.Range("I63").Value = 6050
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 18150
$ws.Range("L63").Value = 30000
$ws.Range("M63").Value = -17401
$ws.Range("N63").Value = -31498

$ws.Range("H66").Value = 6840
$ws.Range("I66").Value = 6050
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 54450
$ws.Range("L66").Value = 90000
$ws.Range("M66").Value = -50706
$ws.Range("N66").Value = -97488

$ws.Range("H114").Value = 8852.5625
$ws.Range("I114").Value = 11181.2
$ws.Range("J114").Value = 4971.5
$ws.Range("K114").Value = 33543.60000000001
$ws.Range("L114").Value = 14914.5
$ws.Range("M114").Value = -30289.60000000001
$ws.Range("N114").Value = -21422.5

$ws.Range("H117").Value = 11757.272
$ws.Range("I117").Value = 809.6667
$ws.Range("J117").Value = 15862.625
$ws.Range("K117").Value = 2429.0001
$ws.Range("L117").Value = 47587.875
$ws.Range("M117").Value = 1012.9999
$ws.Range("N117").Value = -54471.875

$ws.Range("H121").Value = 3601.4524
$ws.Range("I121").Value = 50129.5
$ws.Range("J121").Value = 1275.05
$ws.Range("K121").Value = 150388.5
$ws.Range("L121").Value = 3825.15
$ws.Range("M121").Value = -149078.5
$ws.Range("N121").Value = -6445.15

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 426.7
$ws.Range("I107").Value = 383
$ws.Range("J107").Value = 601.5
$ws.Range("K107").Value = 383
$ws.Range("L107").Value = 601.5
$ws.Range("M107").Value = 1537
$ws.Range("N107").Value = -4441.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 867.6087
$ws.Range("I82").Value = 938.1539
$ws.Range("J82").Value = 775.9
$ws.Range("K82").Value = 938.1539
$ws.Range("L82").Value = 775.9
$ws.Range("M82").Value = -577.1539
$ws.Range("N82").Value = -1497.9

$ws.Range("H85").Value = 867.6087
$ws.Range("I85").Value = 938.1539
$ws.Range("J85").Value = 775.9
$ws.Range("K85").Value = 938.1539
$ws.Range("L85").Value = 775.9
$ws.Range("M85").Value = 309.8461
$ws.Range("N85").Value = -3271.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2450
$ws.Range("I62").Value = 2537.5
$ws.Range("J62").Value = 2275
$ws.Range("K62").Value = 2537.5
$ws.Range("L62").Value = 2275
$ws.Range("M62").Value = -1913.5
$ws.Range("N62").Value = -3523

$ws.Range("H65").Value = 2450
$ws.Range("I65").Value = 2537.5
$ws.Range("J65").Value = 2275
$ws.Range("K65").Value = 12687.5
$ws.Range("L65").Value = 11375
$ws.Range("M65").Value = -9567.5
$ws.Range("N65").Value = -17615

$ws.Range("H81").Value = 1481.6364
$ws.Range("I81").Value = 1325
$ws.Range("J81").Value = 1899.3334
$ws.Range("K81").Value = 2650
$ws.Range("L81").Value = 3798.6668
$ws.Range("M81").Value = -1589
$ws.Range("N81").Value = -5920.6668

$ws.Range("H84").Value = 1481.6364
$ws.Range("I84").Value = 1325
$ws.Range("J84").Value = 1899.3334
$ws.Range("K84").Value = 13250
$ws.Range("L84").Value = 18993.334
$ws.Range("M84").Value = -7946
$ws.Range("N84").Value = -29601.334

$ws.Range("H96").Value = 125000500
$ws.Range("I96").Value = 250000000
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 250000000
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = -249998627
$ws.Range("N96").Value = -3746

$ws.Range("H136").Value = 3854.3953
$ws.Range("I136").Value = 4894.3
$ws.Range("J136").Value = 1454.6154
$ws.Range("K136").Value = 14682.9
$ws.Range("L136").Value = 4363.8462
$ws.Range("M136").Value = -9580.7142
